# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c... 🚀
# Bump the FHIR StructureDefinition metadata from 5.0.0 -> 6.0.0, refresh the
# publish date, swap the "Contact" rows for a Publisher/Jurisdiction pair, and
# update the root Extension's Short/Definition text on the Elements sheet.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes
# "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row;
# it is removed entirely (rows below shift up one, A1:B21 -> A1:B20).
$meta.Rows.Item(11).Delete()

# ---- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): Short / Definition now describe the actual
# extension instead of the generic "Extension" / "An Extension" placeholders.
$elements.Range("K2").Value = "Shortterm Care Duration"
$elements.Range("L2").Value = "Number of weeks for which the employee is eligible for short-term disability (STD) benefits"
